# Weekly update: insert a new price-report row for Perejil (Agrícola del Norte
# S.A. de Arica) at the top of the data (row 21, right after the header block
# that's already sorted at rows 2-20), pushing all the existing report rows
# down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 21; this shifts old rows 21..54 down to 22..55
# and (matching Excel's native behaviour) the new row inherits the date
# column's number format from the row below it.
$ws.Rows("21:21").Insert()

# Populate the newly inserted row with this week's reading.
$ws.Cells.Item(21, 1).Value = 1
$ws.Cells.Item(21, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(21, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(21, 4).Value = 45175
$ws.Cells.Item(21, 5).Value = 15
$ws.Cells.Item(21, 6).Value = 100112044
$ws.Cells.Item(21, 7).Value = "Perejil"
$ws.Cells.Item(21, 8).Value = "Sin especificar"
$ws.Cells.Item(21, 9).Value = "Primera"
$ws.Cells.Item(21, 10).Value = 380
$ws.Cells.Item(21, 11).Value = 1300
$ws.Cells.Item(21, 12).Value = 1500
$ws.Cells.Item(21, 13).Value = 1395
$ws.Cells.Item(21, 14).Value = "`$/atado 1,5 a 2 kilos"
$ws.Cells.Item(21, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(21, 16).Value = 698
$ws.Cells.Item(21, 17).Value = 2
$ws.Cells.Item(21, 18).Value = "Hortaliza"
